$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting rows 24-31 down to 25-32
$ws.Rows.Item(24).Insert()

$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44837
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 100112026
$ws.Cells.Item(24, 7).Value = "Haba"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 520
$ws.Cells.Item(24, 11).Value = 8000
$ws.Cells.Item(24, 12).Value = 9000
$ws.Cells.Item(24, 13).Value = 8500
$ws.Cells.Item(24, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(24, 16).Value = 340
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
